$wb = $excel.ActiveWorkbook

# --- UserInfromation sheet: move selection, it loses the active/selected tab ---
$wsUser = $wb.Worksheets.Item("UserInfromation")
$wsUser.Range("E1").Select()

# --- Sheet3: new hyperlink text + label in row 7, becomes the active tab ---
$wsSheet3 = $wb.Worksheets.Item("Sheet3")

$wsSheet3.Columns.Item(2).ColumnWidth = 26.8333333

$wsSheet3.Range("B7").Value = "https://youtu.be/uwDmr8zpsaY?si=j6QGZ3_A16GSvycj"
$wsSheet3.Range("B7").WrapText = $true
$wsSheet3.Rows.Item(7).RowHeight = 86.4

$wsSheet3.Range("C7").Value = "Shotcut"

$wsSheet3.Activate()
$wsSheet3.Range("C7").Select()
